# Update the worked-answer text in each populated cell of the single
# table in the document. Cells are addressed positionally (row, column)
# so that the update is unambiguous even though some old/new values are
# shared between different cells.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "332÷3=110, 2" },
    @{ Row = 1;  Col = 2; Text = "707÷8=88, 3" },
    @{ Row = 1;  Col = 3; Text = "998÷3=332, 2" },
    @{ Row = 1;  Col = 4; Text = "118÷8=14, 6" },
    @{ Row = 1;  Col = 5; Text = "640÷6=106, 4" },

    @{ Row = 5;  Col = 1; Text = "937÷9=104, 1" },
    @{ Row = 5;  Col = 2; Text = "366÷2=183, 0" },
    @{ Row = 5;  Col = 3; Text = "908÷4=227, 0" },
    @{ Row = 5;  Col = 4; Text = "109÷2=54, 1" },
    @{ Row = 5;  Col = 5; Text = "920÷5=184, 0" },

    @{ Row = 9;  Col = 1; Text = "132÷6=22, 0" },
    @{ Row = 9;  Col = 2; Text = "198÷5=39, 3" },
    @{ Row = 9;  Col = 3; Text = "646÷6=107, 4" },
    @{ Row = 9;  Col = 4; Text = "681÷7=97, 2" },
    @{ Row = 9;  Col = 5; Text = "601÷4=150, 1" },

    @{ Row = 13; Col = 1; Text = "622÷9=69, 1" },
    @{ Row = 13; Col = 2; Text = "706÷5=141, 1" },
    @{ Row = 13; Col = 3; Text = "765÷3=255, 0" },
    @{ Row = 13; Col = 4; Text = "254÷9=28, 2" },
    @{ Row = 13; Col = 5; Text = "912÷6=152, 0" },

    @{ Row = 17; Col = 1; Text = "103÷2=51, 1" },
    @{ Row = 17; Col = 2; Text = "310÷4=77, 2" },
    @{ Row = 17; Col = 3; Text = "308÷4=77, 0" },
    @{ Row = 17; Col = 4; Text = "931÷5=186, 1" },
    @{ Row = 17; Col = 5; Text = "214÷2=107, 0" }
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark/paragraph-mark characters so only the
    # visible text is replaced, keeping the run/paragraph formatting intact.
    $rng.End = $rng.End - 1
    $rng.Text = $u.Text
}
